# Corrected excel sheets for application fix issues
#
# - Split the "Edit Repayment Schedule" sheet: the Approve/Disburse rows
#   (11:16) move out into their own new "Sheet1" tab, which becomes the
#   active tab.

$wb = $excel.ActiveWorkbook

$editSheet = $wb.Worksheets.Item("Edit Repayment Schedule")

# New sheet goes right after "Edit Repayment Schedule" (i.e. at the end).
$newSheet = $wb.Worksheets.Add($null, $editSheet)

# Move (cut/paste) the Approve/Disburse block (A11:B16) onto the new sheet
# as A1:B6, preserving values + formatting.
$editSheet.Range("A11:B16").Cut($newSheet.Range("A1:B6")) | Out-Null

# Remove the now-empty rows from the source sheet so it shrinks back down
# to A1:B10.
$editSheet.Rows("11:16").Delete() | Out-Null

# Restore/update each sheet's own selection.
$editSheet.Activate() | Out-Null
$editSheet.Range("A10:B10").Select() | Out-Null

# New sheet ends up the active tab, with a single-cell selection.
$newSheet.Activate() | Out-Null
$newSheet.Range("A5").Select() | Out-Null
